$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 2027754966
$ws.Range("D2").Value = 581338128

$ws.Range("C3").Value = 8186693023
$ws.Range("D3").Value = 2077596796

$ws.Range("C4").Value = 4768579997
$ws.Range("D4").Value = 1138559903

$ws.Range("C5").Value = 7170315489
$ws.Range("D5").Value = 1590342610

$ws.Range("C6").Value = 8122928131
$ws.Range("D6").Value = 1798640824

$excel.Calculate()
